{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block, plus the\n// blank paragraph that preceded it, right after the\n// \"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito)\" paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the index of the \"LOB1008...\" requirement paragraph so we can also\n// remove the single blank paragraph that directly follows it (part of the\n// same footer block being dropped).\nlet reqIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito)\") !== -1) {\n    reqIdx = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (reqIdx !== -1 && reqIdx + 1 < items.length && items[reqIdx + 1].text === \"\") {\n  toDelete.push(items[reqIdx + 1]);\n}\nfor (let i = 0; i < items.length; i++) {\n  if (targetTexts.indexOf(items[i].text) !== -1) {\n    toDelete.push(items[i]);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the site-footer block that trails the \"Requisitos\" section:\n#   - the blank paragraph right after \"LOB1008: Ci\u00eancia, Tecnologia e\n#     Sociedade (Requisito)\"\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#     pages. Original theme under Creative Commons Attribution\"\n$d = $word.ActiveDocument\n\n$reqMarker = \"LOB1008: Ci\u00eancia, Tecnologia e Sociedade (Requisito)\"\n$copyright = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n$targets = @(\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  $copyright\n)\n\n$n = $d.Paragraphs.Count\n\n# Locate the \"LOB1008...\" paragraph so the blank paragraph immediately\n# following it (part of the removed footer block) can be identified too.\n$reqIndex = -1\nfor ($i = 1; $i -le $n; $i++) {\n  $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $reqMarker) {\n    $reqIndex = $i\n    break\n  }\n}\n\n# Walk backwards so earlier deletions never invalidate later indices.\nfor ($i = $n; $i -ge 1; $i--) {\n  $p = $d.Paragraphs($i)\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  $isBlankAfterReq = ($reqIndex -gt 0) -and ($i -eq ($reqIndex + 1)) -and ($t -eq \"\")\n  if (($targets -contains $t) -or $isBlankAfterReq) {\n    $p.Range.Delete()\n  }\n}\n"}
